# Update the "F" column (想去人数 / "want to go" counter) values that changed
# between the previous scrape and the newly generated output (commit:
# "Update gh-pages to output generated at 456a3b4").
#
# Sheet 1 = 展览 (worksheet "sheet1.xml")
# Sheet 4 = 全部类型 (worksheet "sheet4.xml")
# Sheets 2/3 are untouched by this update.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 2235
$ws1.Range("F4").Value  = 13502
$ws1.Range("F7").Value  = 530
$ws1.Range("F8").Value  = 494
$ws1.Range("F10").Value = 1015
$ws1.Range("F11").Value = 13846
$ws1.Range("F12").Value = 14593
$ws1.Range("F13").Value = 44
$ws1.Range("F16").Value = 45
$ws1.Range("F25").Value = 61
$ws1.Range("F26").Value = 5600
$ws1.Range("F27").Value = 942
$ws1.Range("F28").Value = 1042
$ws1.Range("F29").Value = 5367
$ws1.Range("F31").Value = 31
$ws1.Range("F32").Value = 181

# ---- Sheet 4: 全部类型 ----
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 2235
$ws4.Range("F4").Value  = 13502
$ws4.Range("F8").Value  = 530
$ws4.Range("F9").Value  = 494
$ws4.Range("F11").Value = 1015
$ws4.Range("F12").Value = 13846
$ws4.Range("F13").Value = 14593
$ws4.Range("F14").Value = 44
$ws4.Range("F17").Value = 45
$ws4.Range("F26").Value = 61
$ws4.Range("F27").Value = 5600
$ws4.Range("F28").Value = 942
$ws4.Range("F29").Value = 1042
$ws4.Range("F30").Value = 5367
$ws4.Range("F32").Value = 31
$ws4.Range("F33").Value = 181
